$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cases_analysis")

# Clear the old "Times / Operation conditions / Q" header block (I1:L1)
$ws.Range("I1:L1").ClearContents()

# New compact header + values for Times / Operation conditions / Q, moved to Q1:S2
$ws.Range("Q1").Value = "Times"
$ws.Range("R1").Value = "Operation conditions "
$ws.Range("S1").Value = "Q"
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 20
$ws.Range("S2").Value = 10

# New per-network column headers (F2:I2)
$ws.Range("F2").Value = "S-SN"
$ws.Range("G2").Value = "D-SN"
$ws.Range("H2").Value = "S-NC"
$ws.Range("I2").Value = "D-NC"

# New per-network data values
$ws.Range("F3").Value = 0.27
$ws.Range("G3").Value = 1.56
$ws.Range("H3").Value = 1.64
$ws.Range("I3").Value = 12.85

$ws.Range("F4").Value = 0.6
$ws.Range("G4").Value = 3.31
$ws.Range("H4").Value = 3.57
$ws.Range("I4").Value = 8.86

$ws.Range("F5").Value = 3.26
$ws.Range("G5").Value = 19.32
$ws.Range("H5").Value = 43.21
$ws.Range("I5").Value = 333.39

# J5 held the old, now-displaced value - clear it
$ws.Range("J5").ClearContents()

$ws.Range("K4").Select()
